$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: add the "EXP" column value (O) next to the existing Pin Header entry row
$ws.Range("C10").Value = "O"

# Row 11: new component row - Pin Header / PIN library / devicemart link
$ws.Range("A11").Value = "Pin Header"
$ws.Range("B11").Value = "PIN"
$ws.Range("D11").Value = "https://www.devicemart.co.kr/goods/view?no=5810"

# Move the active selection to D11 (matches the saved cursor position)
[void]$ws.Range("D11").Select()
